# Generate Report for Handback
# Update timestamps / status for the d47dd2d9 and d79a76ad rows
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-07 16:25:38"
$wsOverview.Range("G4").Value = "2016-09-07 16:25:38"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "mt"
$wsZh.Range("E4").Value = "mt"
$wsZh.Range("H3").Value = "2016-09-07 16:25:32"
$wsZh.Range("H4").Value = "2016-09-07 16:25:32"
$wsZh.Range("K3").Value = "2016-09-07 16:25:51"
$wsZh.Range("K4").Value = "2016-09-07 16:25:51"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "mt"
$wsDe.Range("E4").Value = "mt"
$wsDe.Range("H3").Value = "2016-09-07 16:25:38"
$wsDe.Range("H4").Value = "2016-09-07 16:25:38"
$wsDe.Range("K3").Value = "2016-09-07 16:25:59"
$wsDe.Range("K4").Value = "2016-09-07 16:25:59"
